# PAP.xlsx edit: update ticket statuses, add "NR" (Não Resolver) legend entry,
# and move the active selection to K6 (matching the commit's recorded state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# --- Update status (column K) values for several rows ---
$ws.Range("K6").Value  = "EP"   # Ligação dinamica á base de dados  : C  -> EP
$ws.Range("K9").Value  = "C"    # Obrigatorio ter um contacto        : PR -> C
$ws.Range("K10").Value = "C"    # Erro a inserir caso não tenha cod postal : PR -> C
$ws.Range("K13").Value = "C"    # Ao limpar o codigo postal, nem tudo limpa : PR -> C
$ws.Range("K14").Value = "C"    # Funcionamento estranho do codigo postal : EP -> C
$ws.Range("K17").Value = "NR"   # Indicador do telefone com o Codigo Postal : PR -> NR
$ws.Range("K18").Value = "C"    # Voltar ao menu apaga os dados introduzidos : PR -> C

# --- Add a new "NR" / "Não Resolver" row to the status legend (row 13) ---
$ws.Range("A13").Value = "NR"
$ws.Range("B13").Value = "Não Resolver"
$ws.Range("B13:D13").Merge()
$ws.Range("B13:D13").HorizontalAlignment = -4108

# --- Update the active selection to match the saved workbook state ---
$ws.Range("K6").Select()
